$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.663.56'
$ws.Range("E2").Value = '  +3.22%  '
$ws.Range("D3").Value = '3.401.67'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.84%  '
$ws.Range("D8").Value = '3.393.75'
$ws.Range("E8").Value = '  +1.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +13.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.633'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.10%  '
$ws.Range("E13").Value = '  +5.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.68%  '
$ws.Range("D15").Value = '3.951.04'
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.59%  '
$ws.Range("D17").Value = '3.401.20'
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("D19").Value = '65.630.31'
$ws.Range("E19").Value = '  +3.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.992'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '472.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +22.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.22%  '
$ws.Range("E29").Value = '  +4.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.55%  '
$ws.Range("E31").Value = '  +5.13%  '
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '62.86'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '578.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -3.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("D40").Value = '0.0₃0758'
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.374'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("D42").Value = '3.090.24'
$ws.Range("E42").Value = '  -1.69%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("E44").Value = '  +1.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0418'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.97%  '
$ws.Range("E46").Value = '  +3.75%  '
$ws.Range("E47").Value = '  +6.55%  '
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '136.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.08%  '
